# "adding a couple of leagues" - rename "Washington Football Team" to "Washington"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B34").Value = "Washington"
$ws.Range("B35").Value = "Washington"

# Reflect the author's final scroll position / selection as closely as possible.
[void]$ws.Range("B36").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
